$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = "2022/3/9完成"
$ws.Range("B4").Select()
